$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.04027
$ws.Range("H2").Value = 0.12081
$ws.Range("I2").Value = 0.01318991723029425
$ws.Range("J2").Value = 0.01318991723029425
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.007258333333333333
$ws.Range("N2").Value = 0.021775
$ws.Range("O2").Value = 0.000328667160253549
$ws.Range("P2").Value = 0.000328667160253549
$ws.Range("Q2").Value = 0.0002922930833333333
$ws.Range("R2").Value = 0.00263063775
$ws.Range("S2").Value = 0.000004335092640060169
$ws.Range("T2").Value = 0.000004335092640060168

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.04027
$ws.Range("H3").Value = 0.12081
$ws.Range("I3").Value = 0.01318991723029425
$ws.Range("J3").Value = 0.01318991723029425
$ws.Range("O3").Value = 0.7778551418094273
$ws.Range("P3").Value = 0.7778551418094272
$ws.Range("Q3").Value = 0.6917687718200001
$ws.Range("R3").Value = 6.22591894638
$ws.Range("S3").Value = 0.01025984493762514
$ws.Range("T3").Value = 0.01025984493762514

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.04027
$ws.Range("H4").Value = 0.12081
$ws.Range("I4").Value = 0.01318991723029425
$ws.Range("J4").Value = 0.01318991723029425
$ws.Range("M4").Value = 4.898620999999999
$ws.Range("N4").Value = 14.695863
$ws.Range("O4").Value = 0.2218161910303192
$ws.Range("P4").Value = 0.2218161910303192
$ws.Range("Q4").Value = 0.19726746767
$ws.Range("R4").Value = 1.77540720903
$ws.Range("S4").Value = 0.002925737200029049
$ws.Range("T4").Value = 0.002925737200029049

# Row 5
$ws.Range("G5").Value = 0.9943730000000001
$ws.Range("I5").Value = 0.3256940046198011
$ws.Range("J5").Value = 0.325694004619801
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.007258333333333333
$ws.Range("N5").Value = 0.021775
$ws.Range("O5").Value = 0.000328667160253549
$ws.Range("P5").Value = 0.000328667160253549
$ws.Range("Q5").Value = 0.007217490691666667
$ws.Range("R5").Value = 0.06495741622500001
$ws.Range("S5").Value = 0.0001070449236099963
$ws.Range("T5").Value = 0.0001070449236099963

# Row 6
$ws.Range("G6").Value = 0.9943730000000001
$ws.Range("I6").Value = 0.3256940046198011
$ws.Range("J6").Value = 0.325694004619801
$ws.Range("O6").Value = 0.7778551418094273
$ws.Range("P6").Value = 0.7778551418094272
$ws.Range("S6").Value = 0.2533427561500156
$ws.Range("T6").Value = 0.2533427561500156

# Row 7
$ws.Range("G7").Value = 0.9943730000000001
$ws.Range("I7").Value = 0.3256940046198011
$ws.Range("J7").Value = 0.325694004619801
$ws.Range("M7").Value = 4.898620999999999
$ws.Range("N7").Value = 14.695863
$ws.Range("O7").Value = 0.2218161910303192
$ws.Range("P7").Value = 0.2218161910303192
$ws.Range("Q7").Value = 4.871056459632999
$ws.Range("R7").Value = 43.839508136697
$ws.Range("S7").Value = 0.07224420354617547
$ws.Range("T7").Value = 0.07224420354617546

# Row 8
$ws.Range("G8").Value = 2.018446666666666
$ws.Range("H8").Value = 6.055339999999999
$ws.Range("I8").Value = 0.6611160781499047
$ws.Range("J8").Value = 0.6611160781499047
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.007258333333333333
$ws.Range("N8").Value = 0.021775
$ws.Range("O8").Value = 0.000328667160253549
$ws.Range("P8").Value = 0.000328667160253549
$ws.Range("Q8").Value = 0.01465055872222222
$ws.Range("R8").Value = 0.1318550285
$ws.Range("S8").Value = 0.0002172871440034926
$ws.Range("T8").Value = 0.0002172871440034925

# Row 9
$ws.Range("G9").Value = 2.018446666666666
$ws.Range("H9").Value = 6.055339999999999
$ws.Range("I9").Value = 0.6611160781499047
$ws.Range("J9").Value = 0.6611160781499047
$ws.Range("O9").Value = 0.7778551418094273
$ws.Range("P9").Value = 0.7778551418094272
$ws.Range("Q9").Value = 34.67341374681333
$ws.Range("R9").Value = 312.06072372132
$ws.Range("S9").Value = 0.5142525407217865
$ws.Range("T9").Value = 0.5142525407217865

# Row 10
$ws.Range("G10").Value = 2.018446666666666
$ws.Range("H10").Value = 6.055339999999999
$ws.Range("I10").Value = 0.6611160781499047
$ws.Range("J10").Value = 0.6611160781499047
$ws.Range("M10").Value = 4.898620999999999
$ws.Range("N10").Value = 14.695863
$ws.Range("O10").Value = 0.2218161910303192
$ws.Range("P10").Value = 0.2218161910303192
$ws.Range("Q10").Value = 9.88760522871333
$ws.Range("R10").Value = 88.98844705841998
$ws.Range("S10").Value = 0.1466462502841147
$ws.Range("T10").Value = 0.1466462502841147
